# Update cryptos list with latest prices and 1h volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.280.71'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '2.060.11'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.624'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.60%  '
$ws.Range("E8").Value = '  +1.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.383'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.03'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.20%  '
$ws.Range("E11").Value = '  +0.12%  '
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").Value = '2.363.29'
$ws.Range("E13").Value = '  +0.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.76'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.778'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.82%  '
$ws.Range("E17").Value = '  -1.05%  '
$ws.Range("D18").Value = '2.060.01'
$ws.Range("E18").Value = '  +0.87%  '
$ws.Range("D19").Value = '37.205.19'
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.37'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.75%  '
$ws.Range("D22").Value = '0.0₃0810'
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '226.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.70%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.43'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.76%  '
$ws.Range("E28").Value = '  +6.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.77'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("E30").Value = '  -1.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.04'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.117'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.45'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.60%  '
$ws.Range("E34").Value = '  +1.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.60'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.49'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.25'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.72'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.88%  '
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.469.18'
$ws.Range("E42").Value = '  -0.94%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '96.30'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.00%  '
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0935'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.13%  '
$ws.Range("E45").Value = '  +2.82%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.17'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.55%  '
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.27'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.28%  '
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '15.11'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.93%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.16'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.99%  '
$ws.Range("E51").Value = '  +1.11%  '
